# Lemon shark leslie matrix workbook update
# - rename sheet "matrix" -> "LemonSharkLeslieMatrix"
# - add a new sheet "WaplesLifeTable" with a Waples-style life table
# - tweak the selection on the original sheet

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# 1) Rename the original worksheet (defined names auto-update to the new name).
$ws1.Name = "LemonSharkLeslieMatrix"

# 2) Insert the new worksheet right after the renamed one.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "WaplesLifeTable"

# --- populate the WaplesLifeTable sheet -----------------------------------

$ws2.Range("A1").Value = "Age"
$ws2.Range("B1").Value = "sx"
$ws2.Range("C1").Value = "lx"
$ws2.Range("D1").Value = "Nx"
$ws2.Range("E1").Value = "bx"
$ws2.Range("F1").Value = "Nxbx"
$ws2.Range("G1").Value = "ERRO"

# Age column: 1..10
$ws2.Range("A2").Value = 1
$ws2.Range("A3").Value = 2
$ws2.Range("A4").Value = 3
$ws2.Range("A5").Value = 4
$ws2.Range("A6").Value = 5
$ws2.Range("A7").Value = 6
$ws2.Range("A8").Value = 7
$ws2.Range("A9").Value = 8
$ws2.Range("A10").Value = 9
$ws2.Range("A11").Value = 10

# Survival rate sx
$ws2.Range("B2:B10").Value = 0.7
$ws2.Range("B11").Value = 0

# Survivorship lx
$ws2.Range("C2").Value = 1
$ws2.Range("C3").Formula = "=C2*B2"
$ws2.Range("C4:C11").Formula = "=C3*B3"

# Abundance at age Nx
$ws2.Range("D2").Value = 1000
$ws2.Range("D3").Formula = "=D2*B2"
$ws2.Range("D4").Formula = "=D3*B3"
$ws2.Range("D5:D11").Formula = "=D4*B4"

# Births per female bx
$ws2.Range("E2:E3").Value = 0
$ws2.Range("E4:E11").Value = 1

# Nx*bx
$ws2.Range("F2").Formula = "=E2*D2"
$ws2.Range("F3:F11").Formula = "=E3*D3"

# ERRO
$ws2.Range("G4").Formula = "=E4/F`$12"
$ws2.Range("G5:G11").Formula = "=E5/F`$12"

# Totals row
$ws2.Range("D12").Formula = "=SUM(D2:D11)"
$ws2.Range("F12").Formula = "=SUM(F4:F11)"

# Number formats matching the Leslie-matrix sheet's conventions
$ws2.Range("C2:C11").NumberFormat = "0.00"
$ws2.Range("D2:D11").NumberFormat = "0"
$ws2.Range("D12").NumberFormat = "0"
$ws2.Range("F2:F12").NumberFormat = "0"

# Explanatory note for ERRO, wrapped in a wide column
$ws2.Range("H12").Value = "ERRO is the births per female (bx) divided by the total reproductive output. Here, fecundity is constant across age classes, so each age class has the same ERRO."
$ws2.Range("H12").WrapText = $true
$ws2.Columns.Item(8).ColumnWidth = 82.75
$ws2.Rows.Item(12).RowHeight = 31.5

$ws2.Range("C11").Select()

# --- restore the active sheet / selection on the Leslie-matrix sheet ------

$ws1.Range("F10").Select()
$ws1.Activate()
